# Add a new effort log entry (row 15) to the worksheet, mirroring the
# existing rows: Date (col A), Effort [h] (col B), Additional Effort [h]
# (col C), Task (col D).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 15

# Date value (2013-06-18) - stored as the Excel serial date number, same as
# the other date cells in column A (which already carry the date style).
$ws.Cells.Item($row, 1).Value = 41443

# Effort [h]
$ws.Cells.Item($row, 2).Value = 1.5

# Additional Effort [h]
$ws.Cells.Item($row, 3).Value = 2.5

# Task - re-use the same text as the row above (row 14 / "Implementation tc12")
$ws.Cells.Item($row, 4).Value = "Implementation tc12"

# Keep the active selection consistent with the newly added last row
$ws.Range("A15").Select()
